$wb = $excel.ActiveWorkbook

# Rename existing sheet
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Classification Report"

# Add new sheet for confusion matrix, placed right after "Classification Report"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Confusion Matrix"

# Header row
$ws2.Range("B1").Value = "Predicted 1"
$ws2.Range("C1").Value = "Predicted 2"
$ws2.Range("D1").Value = "Predicted 3"

# Row 2
$ws2.Range("A2").Value = "Actual 1"
$ws2.Range("B2").Value = 6524
$ws2.Range("C2").Value = 1
$ws2.Range("D2").Value = 0

# Row 3
$ws2.Range("A3").Value = "Actual 2"
$ws2.Range("B3").Value = 109
$ws2.Range("C3").Value = 7421
$ws2.Range("D3").Value = 53

# Row 4
$ws2.Range("A4").Value = "Actual 3"
$ws2.Range("B4").Value = 0
$ws2.Range("C4").Value = 56
$ws2.Range("D4").Value = 172

# Apply styling (bold, centered, thin border) matching header/label cells
$headerRange = $ws2.Range("B1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

$labelRange = $ws2.Range("A2:A4")
$labelRange.Font.Bold = $true
$labelRange.HorizontalAlignment = -4108
$labelRange.VerticalAlignment = -4160
$labelRange.Borders.LineStyle = 1
$labelRange.Borders.Weight = 2

# Match the page margins used on the rest of the workbook (0.75"/1"/0.5")
$ws2.PageSetup.LeftMargin = 54
$ws2.PageSetup.RightMargin = 54
$ws2.PageSetup.TopMargin = 72
$ws2.PageSetup.BottomMargin = 72
$ws2.PageSetup.HeaderMargin = 36
$ws2.PageSetup.FooterMargin = 36
